$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1740
$ws.Range("I40").Value = 1659.8889
$ws.Range("J40").Value = 1820.1111
$ws.Range("K40").Value = 1659.8889
$ws.Range("L40").Value = 1820.1111
$ws.Range("M40").Value = -1484.8889
$ws.Range("N40").Value = -2170.1111
# Row 98
$ws.Range("H98").Value = 4641.6313
$ws.Range("I98").Value = 2959.9033
$ws.Range("J98").Value = 6646.769
$ws.Range("K98").Value = 2959.9033
$ws.Range("L98").Value = 6646.769
$ws.Range("M98").Value = -1461.9033
$ws.Range("N98").Value = -9642.769
# Row 103
$ws.Range("H103").Value = 7588.4614
$ws.Range("I103").Value = 708.3333
$ws.Range("J103").Value = 13485.714
$ws.Range("K103").Value = 2124.9999
$ws.Range("L103").Value = 40457.142
$ws.Range("M103").Value = -1538.9999
$ws.Range("N103").Value = -41629.142
# Row 112
$ws.Range("H112").Value = 1289.9323
$ws.Range("J112").Value = 1289.9323
$ws.Range("L112").Value = 3869.7969
$ws.Range("N112").Value = -6085.796899999999
# Row 122
$ws.Range("H122").Value = 4641.6313
$ws.Range("I122").Value = 2959.9033
$ws.Range("J122").Value = 6646.769
$ws.Range("K122").Value = 8879.7099
$ws.Range("L122").Value = 19940.307
$ws.Range("M122").Value = -6429.7099
$ws.Range("N122").Value = -24840.307

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 9238195
$ws.Range("I63").Value = 19789632
$ws.Range("J63").Value = 5687.5
$ws.Range("K63").Value = 19789632
$ws.Range("L63").Value = 5687.5
$ws.Range("M63").Value = -19788946
$ws.Range("N63").Value = -7059.5
# Row 66
$ws.Range("H66").Value = 9238195
$ws.Range("I66").Value = 19789632
$ws.Range("J66").Value = 5687.5
$ws.Range("K66").Value = 98948160
$ws.Range("L66").Value = 28437.5
$ws.Range("M66").Value = -98944728
$ws.Range("N66").Value = -35301.5
# Row 74
$ws.Range("H74").Value = 2132.762
$ws.Range("I74").Value = 1397.375
$ws.Range("J74").Value = 4486
$ws.Range("K74").Value = 1397.375
$ws.Range("L74").Value = 4486
$ws.Range("M74").Value = -523.375
$ws.Range("N74").Value = -6234
# Row 77
$ws.Range("H77").Value = 2132.762
$ws.Range("I77").Value = 1397.375
$ws.Range("J77").Value = 4486
$ws.Range("K77").Value = 6986.875
$ws.Range("L77").Value = 22430
$ws.Range("M77").Value = -2618.875
$ws.Range("N77").Value = -31166
# Row 88
$ws.Range("H88").Value = 4764968
$ws.Range("I88").Value = 16668166
$ws.Range("J88").Value = 3688.6
$ws.Range("K88").Value = 16668166
$ws.Range("L88").Value = 3688.6
$ws.Range("M88").Value = -16667760
$ws.Range("N88").Value = -4500.6
# Row 91
$ws.Range("H91").Value = 4764968
$ws.Range("I91").Value = 16668166
$ws.Range("J91").Value = 3688.6
$ws.Range("K91").Value = 16668166
$ws.Range("L91").Value = 3688.6
$ws.Range("M91").Value = -16666762
$ws.Range("N91").Value = -6496.6
# Row 112
$ws.Range("H112").Value = 31216.96
$ws.Range("J112").Value = 31216.96
$ws.Range("L112").Value = 31216.96
$ws.Range("N112").Value = -34170.96
# Row 132
$ws.Range("H132").Value = 2287.56
$ws.Range("I132").Value = 1299.5
$ws.Range("J132").Value = 4828.2856
$ws.Range("K132").Value = 3898.5
$ws.Range("L132").Value = 14484.8568
$ws.Range("M132").Value = -1368.5
$ws.Range("N132").Value = -19544.8568

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3317.087
$ws.Range("I134").Value = 1737.7222
$ws.Range("K134").Value = 5213.1666
$ws.Range("M134").Value = -2678.1666

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 789.1739
$ws.Range("I22").Value = 478.0625
$ws.Range("J22").Value = 1500.2858
$ws.Range("K22").Value = 478.0625
$ws.Range("L22").Value = 1500.2858
$ws.Range("M22").Value = -128.0625
$ws.Range("N22").Value = -2200.2858
# Row 31
$ws.Range("H31").Value = 5112.0615
$ws.Range("I31").Value = 2014.1
$ws.Range("K31").Value = 2014.1
$ws.Range("M31").Value = -1719.1
# Row 34
$ws.Range("H34").Value = 5112.0615
$ws.Range("I34").Value = 2014.1
$ws.Range("K34").Value = 2014.1
$ws.Range("M34").Value = -1812.1
# Row 51
$ws.Range("H51").Value = 23729.25
$ws.Range("J51").Value = 23729.25
$ws.Range("L51").Value = 23729.25
$ws.Range("N51").Value = -25201.25
# Row 59
$ws.Range("H59").Value = 32848.11
$ws.Range("J59").Value = 32848.11
$ws.Range("L59").Value = 32848.11
$ws.Range("N59").Value = -35138.11
# Row 61
$ws.Range("H61").Value = 23729.25
$ws.Range("J61").Value = 23729.25
$ws.Range("L61").Value = 23729.25
$ws.Range("N61").Value = -24425.25
# Row 122
$ws.Range("H122").Value = 2466.9412
$ws.Range("I122").Value = 1479.75
$ws.Range("J122").Value = 3344.4443
$ws.Range("K122").Value = 4439.25
$ws.Range("L122").Value = 10033.3329
$ws.Range("M122").Value = -1989.25
$ws.Range("N122").Value = -14933.3329
# Row 134
$ws.Range("H134").Value = 4612.6177
$ws.Range("I134").Value = 4672.5356
$ws.Range("J134").Value = 4333
$ws.Range("K134").Value = 14017.6068
$ws.Range("L134").Value = 12999
$ws.Range("M134").Value = -11482.6068
$ws.Range("N134").Value = -18069

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 538.75
$ws.Range("I68").Value = 538.75
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1616.25
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -805.25
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 538.75
$ws.Range("I71").Value = 538.75
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 4848.75
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -792.75
$ws.Range("N71").ClearContents()
# Row 113
$ws.Range("H113").Value = 5000566
$ws.Range("I113").Value = 606.3333
$ws.Range("K113").Value = 1818.9999
$ws.Range("M113").Value = 351.0001
# Row 131
$ws.Range("H131").Value = 775.85
$ws.Range("I131").Value = 310
$ws.Range("J131").Value = 827.6111
$ws.Range("K131").Value = 930
$ws.Range("L131").Value = 2482.8333
$ws.Range("M131").Value = 4110
$ws.Range("N131").Value = -12562.8333
# Row 132
$ws.Range("H132").Value = 2127.6206
$ws.Range("I132").Value = 982.5
$ws.Range("J132").Value = 3196.4
$ws.Range("K132").Value = 8842.5
$ws.Range("L132").Value = 28767.6
$ws.Range("M132").Value = -6312.5
$ws.Range("N132").Value = -33827.60000000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 1708.3334
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 15000
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = -212
$ws.Range("N19").Value = -15576
# Row 80
$ws.Range("H80").Value = 14709794
$ws.Range("I80").Value = 41670236
$ws.Range("J80").Value = 4099.909
$ws.Range("K80").Value = 41670236
$ws.Range("L80").Value = 4099.909
$ws.Range("M80").Value = -41669238
$ws.Range("N80").Value = -6095.909
# Row 83
$ws.Range("H83").Value = 14709794
$ws.Range("I83").Value = 41670236
$ws.Range("J83").Value = 4099.909
$ws.Range("K83").Value = 208351180
$ws.Range("L83").Value = 20499.545
$ws.Range("M83").Value = -208346188
$ws.Range("N83").Value = -30483.545
# Row 126
$ws.Range("H126").Value = 3385.33
$ws.Range("I126").Value = 2788.1233
$ws.Range("K126").Value = 8364.369900000002
$ws.Range("M126").Value = -5894.369900000002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4984.9287
$ws.Range("I7").Value = 2027
$ws.Range("K7").Value = 2027
$ws.Range("M7").Value = -1915
# Row 40
$ws.Range("H40").Value = 6450.8696
$ws.Range("I40").Value = 5468.5
$ws.Range("J40").Value = 13000
$ws.Range("K40").Value = 5468.5
$ws.Range("L40").Value = 13000
$ws.Range("M40").Value = -5332.5
$ws.Range("N40").Value = -13272
# Row 126
$ws.Range("H126").Value = 4984.9287
$ws.Range("I126").Value = 2027
$ws.Range("K126").Value = 6081
$ws.Range("M126").Value = -3611
# Row 132
$ws.Range("H132").Value = 3505.7896
$ws.Range("I132").Value = 2723.3225
$ws.Range("J132").Value = 6971
$ws.Range("K132").Value = 8169.967500000001
$ws.Range("L132").Value = 20913
$ws.Range("M132").Value = -5639.967500000001
$ws.Range("N132").Value = -25973

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2184.75
$ws.Range("I126").Value = 989.6
$ws.Range("J126").Value = 5172.625
$ws.Range("K126").Value = 2968.8
$ws.Range("L126").Value = 15517.875
$ws.Range("M126").Value = -498.8000000000002
$ws.Range("N126").Value = -20457.875
# Row 132
$ws.Range("H132").Value = 10754459
$ws.Range("I132").Value = 893.8
$ws.Range("J132").Value = 30306394
$ws.Range("K132").Value = 2681.4
$ws.Range("L132").Value = 90919182
$ws.Range("M132").Value = -151.3999999999996
$ws.Range("N132").Value = -90924242
